$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data row (row 46): date/time/rate, stored as plain text like the
# existing history rows (leading apostrophe keeps "2025-09-28" from being
# auto-converted to a date serial; resetting the style afterwards avoids
# leaving a stray quote-prefix/text style on the cell).
$ws.Cells.Item(46, 1).Value = "'2025-09-28"
$ws.Cells.Item(46, 1).Style = "Normal"
$ws.Cells.Item(46, 2).Value = "21:23:03"
$ws.Cells.Item(46, 3).Value = "1.00 EUR = 1,623.5666"
